$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Taladro"
$ws.Range("B3").Value = "50 kg"
$ws.Range("C3").Value = "basico"
$ws.Range("D3").Value = "1 Taladro"
$ws.Range("E3").Value = "Herramientas electronicas"
$ws.Range("F3").Value = "10x10x10"
$ws.Range("G3").Value = "pendiente"
